$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Statistics")
$ws2 = $wb.Worksheets.Item("Accidents")

# --- Statistics sheet: overwrite rows 2-16 and append rows 17-31 ---
$ws1.Cells.Item(2, 1).Value = "2024-08-27 14:48:53"
$ws1.Cells.Item(2, 2).Value = 33.82234123357149
$ws1.Cells.Item(2, 3).Value = 5
$ws1.Cells.Item(3, 1).Value = "2024-08-27 14:48:55"
$ws1.Cells.Item(3, 2).Value = 34.70162069399007
$ws1.Cells.Item(3, 3).Value = 6
$ws1.Cells.Item(4, 1).Value = "2024-08-27 14:48:57"
$ws1.Cells.Item(4, 2).Value = 37.03816020955322
$ws1.Cells.Item(4, 3).Value = 9
$ws1.Cells.Item(5, 1).Value = "2024-08-27 14:48:59"
$ws1.Cells.Item(5, 2).Value = 39.21010669201421
$ws1.Cells.Item(5, 3).Value = 11
$ws1.Cells.Item(6, 1).Value = "2024-08-27 14:49:01"
$ws1.Cells.Item(6, 2).Value = 38.46490281879766
$ws1.Cells.Item(6, 3).Value = 14
$ws1.Cells.Item(7, 1).Value = "2024-08-27 14:49:03"
$ws1.Cells.Item(7, 2).Value = 39.29688538284288
$ws1.Cells.Item(7, 3).Value = 15
$ws1.Cells.Item(8, 1).Value = "2024-08-27 14:49:05"
$ws1.Cells.Item(8, 2).Value = 39.27772065768913
$ws1.Cells.Item(8, 3).Value = 17
$ws1.Cells.Item(9, 1).Value = "2024-08-27 14:49:07"
$ws1.Cells.Item(9, 2).Value = 36.67309280329896
$ws1.Cells.Item(9, 3).Value = 18
$ws1.Cells.Item(10, 1).Value = "2024-08-27 14:49:09"
$ws1.Cells.Item(10, 2).Value = 28.06036210561786
$ws1.Cells.Item(10, 3).Value = 19
$ws1.Cells.Item(11, 1).Value = "2024-08-27 14:49:11"
$ws1.Cells.Item(11, 2).Value = 33.72263203589913
$ws1.Cells.Item(11, 3).Value = 18
$ws1.Cells.Item(12, 1).Value = "2024-08-27 14:49:14"
$ws1.Cells.Item(12, 2).Value = 33.58718957253446
$ws1.Cells.Item(12, 3).Value = 19
$ws1.Cells.Item(13, 1).Value = "2024-08-27 14:49:16"
$ws1.Cells.Item(13, 2).Value = 28.28250633813379
$ws1.Cells.Item(13, 3).Value = 17
$ws1.Cells.Item(14, 1).Value = "2024-08-27 14:49:18"
$ws1.Cells.Item(14, 2).Value = 31.77988347368915
$ws1.Cells.Item(14, 3).Value = 20
$ws1.Cells.Item(15, 1).Value = "2024-08-27 14:49:20"
$ws1.Cells.Item(15, 2).Value = 29.75150592701652
$ws1.Cells.Item(15, 3).Value = 20
$ws1.Cells.Item(16, 1).Value = "2024-08-27 14:49:22"
$ws1.Cells.Item(16, 2).Value = 29.91431503230635
$ws1.Cells.Item(16, 3).Value = 23
$ws1.Cells.Item(17, 1).Value = "2024-08-27 14:49:24"
$ws1.Cells.Item(17, 2).Value = 27.26311660794591
$ws1.Cells.Item(17, 3).Value = 22
$ws1.Cells.Item(18, 1).Value = "2024-08-27 14:49:26"
$ws1.Cells.Item(18, 2).Value = 27.44778645551681
$ws1.Cells.Item(18, 3).Value = 22
$ws1.Cells.Item(19, 1).Value = "2024-08-27 14:49:28"
$ws1.Cells.Item(19, 2).Value = 21.12725701441377
$ws1.Cells.Item(19, 3).Value = 22
$ws1.Cells.Item(20, 1).Value = "2024-08-27 14:49:30"
$ws1.Cells.Item(20, 2).Value = 21.47993317810211
$ws1.Cells.Item(20, 3).Value = 23
$ws1.Cells.Item(21, 1).Value = "2024-08-27 14:49:32"
$ws1.Cells.Item(21, 2).Value = 24.38316892234656
$ws1.Cells.Item(21, 3).Value = 21
$ws1.Cells.Item(22, 1).Value = "2024-08-27 14:49:34"
$ws1.Cells.Item(22, 2).Value = 20.78166059720699
$ws1.Cells.Item(22, 3).Value = 18
$ws1.Cells.Item(23, 1).Value = "2024-08-27 14:49:36"
$ws1.Cells.Item(23, 2).Value = 24.06209510791275
$ws1.Cells.Item(23, 3).Value = 20
$ws1.Cells.Item(24, 1).Value = "2024-08-27 14:49:38"
$ws1.Cells.Item(24, 2).Value = 24.88620965920845
$ws1.Cells.Item(24, 3).Value = 21
$ws1.Cells.Item(25, 1).Value = "2024-08-27 14:49:40"
$ws1.Cells.Item(25, 2).Value = 21.39336173405988
$ws1.Cells.Item(25, 3).Value = 21
$ws1.Cells.Item(26, 1).Value = "2024-08-27 14:49:42"
$ws1.Cells.Item(26, 2).Value = 18.95274755377982
$ws1.Cells.Item(26, 3).Value = 21
$ws1.Cells.Item(27, 1).Value = "2024-08-27 14:49:44"
$ws1.Cells.Item(27, 2).Value = 24.43699718643903
$ws1.Cells.Item(27, 3).Value = 20
$ws1.Cells.Item(28, 1).Value = "2024-08-27 14:49:46"
$ws1.Cells.Item(28, 2).Value = 26.24975857312631
$ws1.Cells.Item(28, 3).Value = 21
$ws1.Cells.Item(29, 1).Value = "2024-08-27 14:49:48"
$ws1.Cells.Item(29, 2).Value = 25.98849097744464
$ws1.Cells.Item(29, 3).Value = 21
$ws1.Cells.Item(30, 1).Value = "2024-08-27 14:49:50"
$ws1.Cells.Item(30, 2).Value = 25.74623652075309
$ws1.Cells.Item(30, 3).Value = 20
$ws1.Cells.Item(31, 1).Value = "2024-08-27 14:49:52"
$ws1.Cells.Item(31, 2).Value = 26.82238415272914
$ws1.Cells.Item(31, 3).Value = 22

# --- Accidents sheet: append rows 2-7 ---
$ws2.Cells.Item(2, 1).Value = "2024-08-27 14:49:08"
$ws2.Cells.Item(2, 2).Value = "Car and Car"
$ws2.Cells.Item(2, 3).Value = "26.27 and 29.66"
$ws2.Cells.Item(2, 4).Value = 1
$ws2.Cells.Item(3, 1).Value = "2024-08-27 14:49:09"
$ws2.Cells.Item(3, 2).Value = "Car and Car"
$ws2.Cells.Item(3, 3).Value = "29.56 and 0.00"
$ws2.Cells.Item(3, 4).Value = 1
$ws2.Cells.Item(4, 1).Value = "2024-08-27 14:49:27"
$ws2.Cells.Item(4, 2).Value = "Car and Car"
$ws2.Cells.Item(4, 3).Value = "38.18 and 43.76"
$ws2.Cells.Item(4, 4).Value = 2
$ws2.Cells.Item(5, 1).Value = "2024-08-27 14:49:30"
$ws2.Cells.Item(5, 2).Value = "Car and Car"
$ws2.Cells.Item(5, 3).Value = "0.00 and 26.17"
$ws2.Cells.Item(5, 4).Value = 2
$ws2.Cells.Item(6, 1).Value = "2024-08-27 14:49:40"
$ws2.Cells.Item(6, 2).Value = "Car and Car"
$ws2.Cells.Item(6, 3).Value = "31.16 and 41.78"
$ws2.Cells.Item(6, 4).Value = 3
$ws2.Cells.Item(7, 1).Value = "2024-08-27 14:49:41"
$ws2.Cells.Item(7, 2).Value = "Car and Car"
$ws2.Cells.Item(7, 3).Value = "33.89 and 0.00"
$ws2.Cells.Item(7, 4).Value = 3
